$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Compute DGF Population 1 OFF"
# The OpenFisca variable mapping (column B) for two DGF (dotation globale
# de fonctionnement) line items is moved up two rows, from the "total"
# row onto the first sub-component row, effectively turning off the
# mapping that was driving the computation from the totals and attaching
# it to the first fraction instead.

# --- Dotation de solidarite urbaine ---
# Row 130: "... Montant attribution spontanee DSU"
# Row 131: "... Montant progression de la DSU"
# Row 132: "... Montant total reparti"            <- had "dotation_solidarite_urbaine"
$dsuVar = $ws.Range("B132").Text
$ws.Range("B130").Value = $dsuVar
$ws.Range("B132").Value = ""

# --- Dotation de solidarite rurale Bourg-centre ---
# Row 142: "... Montant de la commune eligible"
# Row 143: "... Montant de la garantie de sortie"
# Row 144: "... Montant global reparti"           <- had "premiere_fraction_dotation_solidarite_rurale"
$dsrVar = $ws.Range("B144").Text
$ws.Range("B142").Value = $dsrVar
$ws.Range("B144").Value = ""

# Move the visible selection/cursor to reflect where the edit was made.
$ws.Range("B142").Select() | Out-Null
